$d = $word.ActiveDocument

# Locate the very last paragraph in the document body (the "Mardi :" entry
# under "Semaine S-0 :") so the new "Mercredi :" paragraph can be appended
# right after it, before the closing sectPr.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$endPos = $lastPara.Range.End

# A zero-length Range built from the document (rather than a collapsed
# Paragraph.Range) inserts cleanly at that exact position without
# clobbering the preceding paragraph's content.
$insertionPoint = $d.Range($endPos, $endPos)

$newParagraphXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Mercredi :</w:t></w:r><w:r><w:t xml:space="preserve"> Je crée une page 404 pour ce site avec le jeu </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PacMan</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
